# Auto-committed on 2022/04/14 週四
#
# Change set (per the author's diff):
#   1. DBD!E16  14 -> 10            (每筆扣款限額 field "長度"/length)
#   2. DBD!F16  (new cell) = 2       (同一列新增 "小數"/decimal-places)
#   3. DBD sheet: active selection moves from G13 to G9
#
# (workbook.xml bits such as fileVersion/rupBuild, the xr:revisionPtr
#  session GUID, and the workbookView window-chrome geometry are stamped
#  by the authoring application itself on every save and are not
#  reachable through the Excel object model - they are intentionally
#  left untouched here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 16 data edits.
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 2

# Move the sheet's active selection to G9 (was G13).
$ws.Activate() | Out-Null
$ws.Range("G9").Select() | Out-Null
